$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 336.73334
$ws.Range("I6").Value = 273.07693
$ws.Range("J6").Value = 750.5
$ws.Range("K6").Value = 819.2307900000001
$ws.Range("L6").Value = 2251.5
$ws.Range("M6").Value = -707.2307900000001
$ws.Range("N6").Value = -2475.5
# Row 125
$ws.Range("H125").Value = 8010519.5
$ws.Range("I125").Value = 1726.2
$ws.Range("J125").Value = 12459849
$ws.Range("K125").Value = 15535.8
$ws.Range("L125").Value = 112138641
$ws.Range("M125").Value = -13075.8
$ws.Range("N125").Value = -112143561
# Row 132
$ws.Range("H132").Value = 28891.158
$ws.Range("I132").Value = 29645.082
$ws.Range("J132").Value = 996
$ws.Range("K132").Value = 88935.246
$ws.Range("L132").Value = 2988
$ws.Range("M132").Value = -86405.246
$ws.Range("N132").Value = -8048
# Row 134
$ws.Range("H134").Value = 77777
$ws.Range("J134").Value = 77777
$ws.Range("L134").Value = 77777
$ws.Range("N134").Value = -87917
# Row 135
$ws.Range("H135").Value = 2673.4546
$ws.Range("I135").Value = 2440.8
$ws.Range("J135").Value = 5000
$ws.Range("K135").Value = 21967.2
$ws.Range("L135").Value = 45000
$ws.Range("M135").Value = -19432.2
$ws.Range("N135").Value = -50070
# Row 137
$ws.Range("H137").Value = 41668892
$ws.Range("I137").Value = 55556676
$ws.Range("J137").Value = 5533.8335
$ws.Range("K137").Value = 166670028
$ws.Range("L137").Value = 16601.5005
$ws.Range("M137").Value = -166667478
$ws.Range("N137").Value = -21701.5005
# Row 138
$ws.Range("H138").Value = 2998.6262
$ws.Range("I138").Value = 2564.3333
$ws.Range("J138").Value = 3215.7727
$ws.Range("K138").Value = 7692.999899999999
$ws.Range("L138").Value = 9647.3181
$ws.Range("M138").Value = -2552.999899999999
$ws.Range("N138").Value = -19927.3181

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 20981.254
$ws.Range("I32").Value = 5979.5796
$ws.Range("K32").Value = 5979.5796
$ws.Range("M32").Value = -5692.5796
# Row 45
$ws.Range("H45").Value = 1080.76
$ws.Range("I45").Value = 1066.0435
$ws.Range("K45").Value = 1066.0435
$ws.Range("M45").Value = -689.0435
# Row 132
$ws.Range("H132").Value = 3591.96
$ws.Range("I132").Value = 3169.35
$ws.Range("J132").Value = 5282.4
$ws.Range("K132").Value = 9508.049999999999
$ws.Range("L132").Value = 15847.2
$ws.Range("M132").Value = -6978.049999999999
$ws.Range("N132").Value = -20907.2
# Row 133
$ws.Range("H133").Value = 52333.332
$ws.Range("J133").Value = 52333.332
$ws.Range("L133").Value = 52333.332
$ws.Range("N133").Value = -57393.332
# Row 139
$ws.Range("H139").Value = 50803.75
$ws.Range("J139").Value = 50803.75
$ws.Range("L139").Value = 50803.75
$ws.Range("N139").Value = -61083.75

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 132
$ws.Range("H132").Value = 45474.285
$ws.Range("J132").Value = 45474.285
$ws.Range("L132").Value = 45474.285
$ws.Range("N132").Value = -55594.285
# Row 134
$ws.Range("H134").Value = 1843.317
$ws.Range("I134").Value = 1742.6154
$ws.Range("J134").Value = 3807
$ws.Range("K134").Value = 5227.8462
$ws.Range("L134").Value = 11421
$ws.Range("M134").Value = -2692.8462
$ws.Range("N134").Value = -16491
# Row 137
$ws.Range("H137").Value = 61390
$ws.Range("J137").Value = 61390
$ws.Range("L137").Value = 61390
$ws.Range("N137").Value = -71590
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
# Row 139
$ws.Range("H139").Value = 94999.5
$ws.Range("J139").Value = 94999.5
$ws.Range("L139").Value = 94999.5
$ws.Range("N139").Value = -105279.5
# Row 140
$ws.Range("H140").Value = 73945
$ws.Range("J140").Value = 73945
$ws.Range("L140").Value = 73945
$ws.Range("N140").Value = -84305
# Row 141
$ws.Range("H141").Value = 87246
$ws.Range("J141").Value = 81780
$ws.Range("L141").Value = 81780
$ws.Range("N141").Value = -92140

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 54
$ws.Range("I7").Value = 53.8
$ws.Range("K7").Value = 53.8
$ws.Range("M7").Value = 59.2
# Row 31
$ws.Range("H31").Value = 4419.0483
$ws.Range("I31").Value = 2154.1667
$ws.Range("J31").Value = 5849.5
$ws.Range("K31").Value = 2154.1667
$ws.Range("L31").Value = 5849.5
$ws.Range("M31").Value = -1859.1667
$ws.Range("N31").Value = -6439.5
# Row 34
$ws.Range("H34").Value = 4419.0483
$ws.Range("I34").Value = 2154.1667
$ws.Range("J34").Value = 5849.5
$ws.Range("K34").Value = 2154.1667
$ws.Range("L34").Value = 5849.5
$ws.Range("M34").Value = -1952.1667
$ws.Range("N34").Value = -6253.5
# Row 134
$ws.Range("H134").Value = 27779918
$ws.Range("I134").Value = 32259526
$ws.Range("J134").Value = 21742186
$ws.Range("K134").Value = 96778578
$ws.Range("L134").Value = 65226558
$ws.Range("M134").Value = -96776043
$ws.Range("N134").Value = -65231628

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 172.9
$ws.Range("J98").Value = 146.57143
$ws.Range("L98").Value = 439.71429
$ws.Range("N98").Value = -3435.71429
# Row 113
$ws.Range("H113").Value = 2841716.5
$ws.Range("I113").Value = 482.9091
$ws.Range("J113").Value = 4329981.5
$ws.Range("K113").Value = 1448.7273
$ws.Range("L113").Value = 12989944.5
$ws.Range("M113").Value = 721.2727
$ws.Range("N113").Value = -12994284.5
# Row 125
$ws.Range("H125").Value = 2265
$ws.Range("I125").Value = 2265
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 6795
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -1875
$ws.Range("N125").ClearContents()
# Row 131
$ws.Range("H131").Value = 10754735
$ws.Range("I131").Value = 943.3333
$ws.Range("J131").Value = 11906927
$ws.Range("K131").Value = 2829.9999
$ws.Range("L131").Value = 35720781
$ws.Range("M131").Value = 2210.0001
$ws.Range("N131").Value = -35730861
# Row 134
$ws.Range("H134").Value = 6932.25
$ws.Range("J134").Value = 10673.267
$ws.Range("L134").Value = 32019.801
$ws.Range("N134").Value = -42159.801
# Row 137
$ws.Range("H137").Value = 7219261.5
$ws.Range("J137").Value = 339688.66
$ws.Range("L137").Value = 1019065.98
$ws.Range("N137").Value = -1029265.98
# Row 138
$ws.Range("H138").Value = 964.2857
$ws.Range("I138").Value = 964.2857
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 2892.8571
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 2247.1429
$ws.Range("N138").ClearContents()
# Row 140
$ws.Range("H140").Value = 6767.4443
$ws.Range("I140").Value = 9412.708000000001
$ws.Range("J140").Value = 3744.2856
$ws.Range("K140").Value = 28238.124
$ws.Range("L140").Value = 11232.8568
$ws.Range("M140").Value = -23058.124
$ws.Range("N140").Value = -21592.8568
# Row 141
$ws.Range("H141").Value = 5959.875
$ws.Range("I141").Value = 6735.8
$ws.Range("K141").Value = 20207.4
$ws.Range("M141").Value = -15027.4

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1797.75
$ws.Range("I113").Value = 1795.5
$ws.Range("K113").Value = 1795.5
$ws.Range("M113").Value = 374.5
# Row 123
$ws.Range("H123").Value = 12180
$ws.Range("J123").Value = 12180
$ws.Range("L123").Value = 12180
$ws.Range("N123").Value = -17080
# Row 126
$ws.Range("H126").Value = 3020.465
$ws.Range("I126").Value = 3073.75
$ws.Range("J126").Value = 2988.889
$ws.Range("K126").Value = 9221.25
$ws.Range("L126").Value = 8966.667000000001
$ws.Range("M126").Value = -6751.25
$ws.Range("N126").Value = -13906.667
# Row 132
$ws.Range("H132").Value = 3002.7742
$ws.Range("I132").Value = 2712.5652
$ws.Range("J132").Value = 3837.125
$ws.Range("K132").Value = 8137.6956
$ws.Range("L132").Value = 11511.375
$ws.Range("M132").Value = -5607.6956
$ws.Range("N132").Value = -16571.375
# Row 138
$ws.Range("H138").Value = 80214.25
$ws.Range("J138").Value = 80214.25
$ws.Range("L138").Value = 80214.25
$ws.Range("N138").Value = -90494.25
# Row 139
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280
# Row 140
$ws.Range("H140").Value = 75832
$ws.Range("J140").Value = 75832
$ws.Range("L140").Value = 75832
$ws.Range("N140").Value = -86192

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 2464.8572
$ws.Range("I93").Value = 2063.5
$ws.Range("K93").Value = 2063.5
$ws.Range("M93").Value = -815.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 982.5769
$ws.Range("I122").Value = 917.8333
$ws.Range("J122").Value = 1254.5
$ws.Range("K122").Value = 2753.4999
$ws.Range("L122").Value = 3763.5
$ws.Range("M122").Value = -303.4998999999998
$ws.Range("N122").Value = -8663.5
# Row 132
$ws.Range("H132").Value = 2753.4753
$ws.Range("I132").Value = 2414.6592
$ws.Range("J132").Value = 3630.4119
$ws.Range("K132").Value = 7243.9776
$ws.Range("L132").Value = 10891.2357
$ws.Range("M132").Value = -4713.9776
$ws.Range("N132").Value = -15951.2357
